# Applies the "Skills audit" table edits:
#  - Luke K / "What would you like to work on?": drop the stray _GoBack bookmark
#  - McCaulay: wrap name + the jargon tokens (php, mysql, c++) in spell-check
#    proofErr markers
#  - Daiyaan / Kleanthis: wrap the names in spell-check proofErr markers
#  - Lewis row: set the row height and fill in the three previously-blank
#    cells (incl. moving the _GoBack bookmark onto "Keeping track of time")

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Luke K row (3): remove the _GoBack bookmark from "The report" ---------
$cell = $t.Cell(3, 4)
$xml = '<w:p xmlns:w="' + $w + '"><w:r><w:t>The report</w:t></w:r></w:p>'
$cell.Range.InsertXML($xml)

# --- McCaulay row (5) --------------------------------------------------------
$cell = $t.Cell(5, 1)
$xml = '<w:p xmlns:w="' + $w + '">' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>McCaulay</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '</w:p>'
$cell.Range.InsertXML($xml)

$cell = $t.Cell(5, 2)
$xml = '<w:p xmlns:w="' + $w + '">' +
       '<w:r><w:t xml:space="preserve">Basic android, java, c#, </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/><w:r><w:t>php</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/><w:r><w:t>mysql</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/><w:r><w:t>c++</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t>, unity</w:t></w:r>' +
       '</w:p>'
$cell.Range.InsertXML($xml)

# --- Daiyaan row (6): wrap the name only ------------------------------------
$cell = $t.Cell(6, 1)
$xml = '<w:p xmlns:w="' + $w + '">' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>Daiyaan</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '</w:p>'
$cell.Range.InsertXML($xml)

# --- Lewis row (7): height + fill the three blank cells ---------------------
$row = $t.Rows.Item(7)
$row.Height = 28.65

$cell = $t.Cell(7, 2)
$xml = '<w:p xmlns:w="' + $w + '">' +
       '<w:r><w:t xml:space="preserve">Java, html, CSS, basic </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/><w:r><w:t>php</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
       '</w:p>'
$cell.Range.InsertXML($xml)

$cell = $t.Cell(7, 3)
$xml = '<w:p xmlns:w="' + $w + '">' +
       '<w:r><w:t>Keeping track of time</w:t></w:r>' +
       '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
       '<w:bookmarkEnd w:id="0"/>' +
       '</w:p>'
$cell.Range.InsertXML($xml)

$cell = $t.Cell(7, 4)
$xml = '<w:p xmlns:w="' + $w + '"><w:r><w:t>Anything</w:t></w:r></w:p>'
$cell.Range.InsertXML($xml)

# --- Kleanthis row (8): wrap the name only -----------------------------------
$cell = $t.Cell(8, 1)
$xml = '<w:p xmlns:w="' + $w + '">' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>Kleanthis</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '</w:p>'
$cell.Range.InsertXML($xml)
